$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in columns A, Q, R between row 3 and row 4
$a3 = $ws.Range("A3").Value2
$a4 = $ws.Range("A4").Value2
$ws.Range("A3").Value2 = $a4
$ws.Range("A4").Value2 = $a3

$q3 = $ws.Range("Q3").Value2
$q4 = $ws.Range("Q4").Value2
$ws.Range("Q3").Value2 = $q4
$ws.Range("Q4").Value2 = $q3

$r3 = $ws.Range("R3").Value2
$r4 = $ws.Range("R4").Value2
$ws.Range("R3").Value2 = $r4
$ws.Range("R4").Value2 = $r3
